# "1st changes of mifos to finflux"
#
# The "Repayment schedule" sheet gains a new (blank) column inserted
# immediately before the existing "Late" column - "Late" shifts from N to O
# and "Outstanding" shifts from P to Q. The sheet becomes the active sheet
# of the workbook, with the new empty column's cell selected.

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N - shifts "Late"/"Outstanding" one
# column to the right (N->O, P->Q) and leaves a blank column behind at N.
$wsSchedule.Columns("N").Insert()

# Make "Repayment schedule" the active sheet/tab, with the new column's
# top data cell selected.
$wsSchedule.Activate()
$wsSchedule.Range("O8").Select()

$wb.Save()
